$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking text (e.g. prices, hour codes).
# They must be written while forcing a Text number format so Excel keeps them
# as text (matching the original inlineStr cells) instead of coercing to Number.
$numericTextEdits = @(
    @{Row=2; Col='D'; Value='238.23'},
    @{Row=2; Col='G'; Value='13'},
    @{Row=3; Col='D'; Value='21.72'},
    @{Row=3; Col='G'; Value='13'},
    @{Row=4; Col='D'; Value='5.445'},
    @{Row=4; Col='G'; Value='13'},
    @{Row=5; Col='D'; Value='0.05643'},
    @{Row=5; Col='G'; Value='13'},
    @{Row=6; Col='G'; Value='13'},
    @{Row=7; Col='D'; Value='3.352'},
    @{Row=7; Col='G'; Value='13'},
    @{Row=8; Col='D'; Value='0.7944'},
    @{Row=8; Col='G'; Value='13'},
    @{Row=9; Col='D'; Value='1.033'},
    @{Row=9; Col='G'; Value='13'},
    @{Row=10; Col='D'; Value='0.1396'},
    @{Row=10; Col='G'; Value='13'},
    @{Row=11; Col='D'; Value='0.07347'},
    @{Row=11; Col='G'; Value='13'},
    @{Row=12; Col='D'; Value='0.03205'},
    @{Row=12; Col='G'; Value='13'},
    @{Row=13; Col='D'; Value='0.02973'},
    @{Row=13; Col='G'; Value='13'},
    @{Row=14; Col='D'; Value='0.09240'},
    @{Row=14; Col='G'; Value='13'},
    @{Row=15; Col='D'; Value='0.001672'},
    @{Row=15; Col='G'; Value='13'},
    @{Row=16; Col='D'; Value='3.261'},
    @{Row=16; Col='G'; Value='13'},
    @{Row=17; Col='D'; Value='0.04774'},
    @{Row=17; Col='G'; Value='13'},
    @{Row=18; Col='D'; Value='0.0005744'},
    @{Row=18; Col='G'; Value='13'},
    @{Row=19; Col='D'; Value='0.006225'},
    @{Row=19; Col='G'; Value='13'},
    @{Row=20; Col='D'; Value='0.005096'},
    @{Row=20; Col='G'; Value='13'},
    @{Row=21; Col='D'; Value='0.001052'},
    @{Row=21; Col='G'; Value='13'},
    @{Row=22; Col='D'; Value='0.0001501'},
    @{Row=22; Col='G'; Value='13'},
    @{Row=23; Col='D'; Value='0.0003212'},
    @{Row=23; Col='G'; Value='13'},
    @{Row=24; Col='D'; Value='3.909'},
    @{Row=24; Col='G'; Value='13'},
    @{Row=25; Col='D'; Value='2.201'},
    @{Row=25; Col='G'; Value='13'},
    @{Row=26; Col='D'; Value='0.3337'},
    @{Row=26; Col='G'; Value='13'},
    @{Row=27; Col='D'; Value='0.1054'},
    @{Row=27; Col='G'; Value='13'},
    @{Row=28; Col='G'; Value='13'},
    @{Row=29; Col='G'; Value='13'},
    @{Row=30; Col='G'; Value='13'},
    @{Row=31; Col='G'; Value='13'},
    @{Row=32; Col='G'; Value='13'},
    @{Row=33; Col='G'; Value='13'},
    @{Row=34; Col='G'; Value='13'},
    @{Row=35; Col='G'; Value='13'},
    @{Row=36; Col='G'; Value='13'},
    @{Row=37; Col='G'; Value='13'},
    @{Row=38; Col='G'; Value='13'},
    @{Row=39; Col='G'; Value='13'},
    @{Row=40; Col='D'; Value='0.04122'},
    @{Row=40; Col='G'; Value='13'},
    @{Row=41; Col='D'; Value='0.006938'},
    @{Row=41; Col='G'; Value='13'},
    @{Row=42; Col='D'; Value='0.1041'},
    @{Row=42; Col='G'; Value='13'},
    @{Row=43; Col='D'; Value='0.003012'},
    @{Row=43; Col='G'; Value='13'},
    @{Row=44; Col='D'; Value='0.009387'},
    @{Row=44; Col='G'; Value='13'},
    @{Row=45; Col='D'; Value='0.00005439'},
    @{Row=45; Col='G'; Value='13'},
    @{Row=46; Col='G'; Value='13'},
    @{Row=47; Col='D'; Value='0.6756'},
    @{Row=47; Col='G'; Value='13'},
    @{Row=48; Col='D'; Value='0.03587'},
    @{Row=48; Col='G'; Value='13'},
    @{Row=49; Col='D'; Value='0.00002101'},
    @{Row=49; Col='G'; Value='13'},
    @{Row=50; Col='G'; Value='13'},
    @{Row=51; Col='G'; Value='13'}
)

foreach ($edit in $numericTextEdits) {
    $cell = $ws.Range($edit.Col + $edit.Row)
    $cell.NumberFormat = "@"
    $cell.Value = $edit.Value
    $cell.ClearFormats()
}

# Cells whose new values are plain text (coin names, URLs, composite labels);
# Excel keeps these as text automatically.
$textEdits = @(
    @{Row=14; Col='B'; Value='BitMartToken'},
    @{Row=14; Col='C'; Value='https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'},
    @{Row=14; Col='E'; Value='13BitMartTokenBMX'},
    @{Row=15; Col='B'; Value='BitForexToken'},
    @{Row=15; Col='C'; Value='https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'},
    @{Row=15; Col='E'; Value='14BitForexTokenBF'},
    @{Row=16; Col='B'; Value='MCDex'},
    @{Row=16; Col='C'; Value='https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'},
    @{Row=16; Col='E'; Value='15MCDexMCB'},
    @{Row=17; Col='B'; Value='CoinExToken'},
    @{Row=17; Col='C'; Value='https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'},
    @{Row=17; Col='E'; Value='16CoinExTokenCET'},
    @{Row=18; Col='B'; Value='One'},
    @{Row=18; Col='C'; Value='https://coinranking.com/coin/6Lga5NiXX3rT+one-one'},
    @{Row=18; Col='E'; Value='17OneONE'},
    @{Row=19; Col='B'; Value='TigerCash'},
    @{Row=19; Col='C'; Value='https://coinranking.com/coin/6hIn06L2+tigercash-tch'},
    @{Row=19; Col='E'; Value='18TigerCashTCH'},
    @{Row=20; Col='B'; Value='HotbitToken'},
    @{Row=20; Col='C'; Value='https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'},
    @{Row=20; Col='E'; Value='19HotbitTokenHTB'},
    @{Row=21; Col='B'; Value='BitKan'},
    @{Row=21; Col='C'; Value='https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'},
    @{Row=21; Col='E'; Value='20BitKanKAN'},
    @{Row=22; Col='B'; Value='NitroEx'},
    @{Row=22; Col='C'; Value='https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'},
    @{Row=22; Col='E'; Value='21NitroExNTX'},
    @{Row=23; Col='B'; Value='UpBots'},
    @{Row=23; Col='C'; Value='https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'},
    @{Row=23; Col='E'; Value='22UpBotsUBXT'},
    @{Row=24; Col='B'; Value='LEO'},
    @{Row=24; Col='C'; Value='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'},
    @{Row=24; Col='E'; Value='23LEOLEOBestin24h'},
    @{Row=25; Col='B'; Value='BTSEToken'},
    @{Row=25; Col='C'; Value='https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'},
    @{Row=25; Col='E'; Value='24BTSETokenBTSE'},
    @{Row=26; Col='B'; Value='BitpandaEcosystemToken'},
    @{Row=26; Col='C'; Value='https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'},
    @{Row=26; Col='E'; Value='25BitpandaEcosystemTokenBEST'},
    @{Row=27; Col='B'; Value='ProBitToken'},
    @{Row=27; Col='C'; Value='https://coinranking.com/coin/lQP4d6T2+probittoken-prob'},
    @{Row=27; Col='E'; Value='26ProBitTokenPROB'},
    @{Row=42; Col='B'; Value='BKEXToken'},
    @{Row=42; Col='C'; Value='https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'},
    @{Row=42; Col='E'; Value='41BKEXTokenBKK'},
    @{Row=43; Col='B'; Value='CEJI'},
    @{Row=43; Col='C'; Value='https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'},
    @{Row=43; Col='E'; Value='42CEJICEJI'},
    @{Row=48; Col='E'; Value='47BOLOBOLOWorstin24h'}
)

foreach ($edit in $textEdits) {
    $ws.Range($edit.Col + $edit.Row).Value = $edit.Value
}

